# Atualização de bases das ligas, do dia: 09-04-2024 às 22:40
#
# The underlying data row-pairs got swapped (their full record, columns B:AC,
# while keeping the sequential index in column A unchanged). Swap back the
# content of each pair so the workbook matches the updated source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (id) through AC (PL_AhUnder) hold the per-match record; column A
# is just a running sequence number and must stay where it is.
$firstCol = 2   # B
$lastCol  = 29  # AC

$rowPairs = @(
    @(5, 6),
    @(9, 10),
    @(16, 17),
    @(25, 26),
    @(37, 38),
    @(63, 64),
    @(88, 89),
    @(92, 93),
    @(99, 100),
    @(103, 104)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($r1, $firstCol), $ws.Cells.Item($r1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($r2, $firstCol), $ws.Cells.Item($r2, $lastCol))

    # .Value's getter round-trips oddly in this host; Value2 behaves like a
    # normal COM Variant (returns a 2-D System.Object[,] for a multi-cell
    # range) and its setter accepts that array straight back, so use it for
    # both sides of the swap.
    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
